# Commit: "wrapping up test file audit"
#
# Content edit recovered from the diff:
#   - In the "optimization_parameters" worksheet, the row holding the
#     leftover "Sheet" / 3 / 4 cells (row 16) is deleted entirely, shifting
#     the "simulation_timepoints" row (old row 17) up to row 16. Deleting
#     this row also makes the shared string "Sheet" and the one-off number
#     style that was only used by that row's B cell become unused, which is
#     why they disappear from sharedStrings.xml / styles.xml on save.
#   - The active / selected sheet moves from "optimization_parameters" to
#     the last sheet, "optimization_diagnostics" (bookViews activeTab moves
#     from 6 to 13), so tabSelected hops from sheet7 to sheet14.

$wb = $excel.ActiveWorkbook

$paramsSheet = $wb.Worksheets.Item("optimization_parameters")

# Delete the whole "Sheet" row (row 16): A16="Sheet", B16=3, C16=4.
$paramsSheet.Rows.Item(16).Delete()

# Keep the selection in sync with the row that shifted into position 16
# (this mirrors what Excel does automatically when a row above the
# selection is removed).
[void]$paramsSheet.Range("A16:XFD16").Select()

# Move focus to the last sheet, which becomes the active tab on save.
[void]$wb.Worksheets.Item("optimization_diagnostics").Activate()
